$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$rng = $ws.Range("D2:G2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = '319.48'
$ws.Range("E2").Value = '-3.81%'
$ws.Range("G2").Value = '5'
$rng.Style = "Normal"

# Row 3
$rng = $ws.Range("D3:G3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = '42.70'
$ws.Range("E3").Value = '-7.19%'
$ws.Range("G3").Value = '5'
$rng.Style = "Normal"

# Row 4
$rng = $ws.Range("D4:G4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = '5.163'
$ws.Range("E4").Value = '-9.09%'
$ws.Range("G4").Value = '5'
$rng.Style = "Normal"

# Row 5
$rng = $ws.Range("D5:G5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = '0.08140'
$ws.Range("E5").Value = '-2.81%'
$ws.Range("G5").Value = '5'
$rng.Style = "Normal"

# Row 6
$rng = $ws.Range("D6:G6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = '4.329'
$ws.Range("E6").Value = '-3.30%'
$ws.Range("G6").Value = '5'
$rng.Style = "Normal"

# Row 7
$rng = $ws.Range("D7:G7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = '1.753'
$ws.Range("E7").Value = '-14.04%'
$ws.Range("G7").Value = '5'
$rng.Style = "Normal"

# Row 8
$rng = $ws.Range("D8:G8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = '0.9514'
$ws.Range("E8").Value = '-3.77%'
$ws.Range("G8").Value = '5'
$rng.Style = "Normal"

# Row 9
$rng = $ws.Range("D9:G9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = '0.1110'
$ws.Range("E9").Value = '-4.80%'
$ws.Range("G9").Value = '5'
$rng.Style = "Normal"

# Row 10
$rng = $ws.Range("D10:G10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = '0.1846'
$ws.Range("E10").Value = '-4.71%'
$ws.Range("G10").Value = '5'
$rng.Style = "Normal"

# Row 11
$rng = $ws.Range("D11:G11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = '0.04667'
$ws.Range("E11").Value = '-0.31%'
$ws.Range("G11").Value = '5'
$rng.Style = "Normal"

# Row 12
$rng = $ws.Range("D12:G12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = '0.09348'
$ws.Range("E12").Value = '-6.25%'
$ws.Range("G12").Value = '5'
$rng.Style = "Normal"

# Row 13
$rng = $ws.Range("D13:G13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = '7.413'
$ws.Range("E13").Value = '-28.59%'
$ws.Range("G13").Value = '5'
$rng.Style = "Normal"

# Row 14
$rng = $ws.Range("D14:G14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = '0.1058'
$ws.Range("E14").Value = '-0.12%'
$ws.Range("G14").Value = '5'
$rng.Style = "Normal"

# Row 15
$rng = $ws.Range("D15:G15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = '0.001283'
$ws.Range("E15").Value = '-0.17%'
$ws.Range("G15").Value = '5'
$rng.Style = "Normal"

# Row 16
$rng = $ws.Range("D16:G16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = '0.005874'
$ws.Range("E16").Value = '-3.66%'
$ws.Range("G16").Value = '5'
$rng.Style = "Normal"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

# Row 17
$rng = $ws.Range("D17:G17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = '0.004328'
$ws.Range("E17").Value = '-6.75%'
$ws.Range("G17").Value = '5'
$rng.Style = "Normal"
$ws.Range("B17").Value = 'HotbitToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'

# Row 18
$rng = $ws.Range("D18:G18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = '3.364'
$ws.Range("E18").Value = '-0.17%'
$ws.Range("G18").Value = '5'
$rng.Style = "Normal"

# Row 19
$rng = $ws.Range("D19:G19")
$rng.NumberFormat = "@"
$ws.Range("G19").Value = '5'
$rng.Style = "Normal"

# Row 20
$rng = $ws.Range("D20:G20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = '0.3364'
$ws.Range("E20").Value = '-0.03%'
$ws.Range("G20").Value = '5'
$rng.Style = "Normal"

# Row 21
$rng = $ws.Range("D21:G21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = '0.1381'
$ws.Range("E21").Value = '-1.36%'
$ws.Range("G21").Value = '5'
$rng.Style = "Normal"

# Row 22
$rng = $ws.Range("D22:G22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = '0.2628'
$ws.Range("E22").Value = '-0.89%'
$ws.Range("G22").Value = '5'
$rng.Style = "Normal"

# Row 23
$rng = $ws.Range("D23:G23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = '0.04187'
$ws.Range("E23").Value = '-0.50%'
$ws.Range("G23").Value = '5'
$rng.Style = "Normal"
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

# Row 24
$rng = $ws.Range("D24:G24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = '0.001252'
$ws.Range("E24").Value = '-4.29%'
$ws.Range("G24").Value = '5'
$rng.Style = "Normal"
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'

# Row 25
$rng = $ws.Range("D25:G25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = '0.0001115'
$ws.Range("E25").Value = '-13.06%'
$ws.Range("G25").Value = '5'
$rng.Style = "Normal"

# Row 26
$rng = $ws.Range("D26:G26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = '0.0002984'
$ws.Range("E26").Value = '-20.33%'
$ws.Range("G26").Value = '5'
$rng.Style = "Normal"

# Row 27
$rng = $ws.Range("D27:G27")
$rng.NumberFormat = "@"
$ws.Range("G27").Value = '5'
$rng.Style = "Normal"

# Row 28
$rng = $ws.Range("D28:G28")
$rng.NumberFormat = "@"
$ws.Range("G28").Value = '5'
$rng.Style = "Normal"

# Row 29
$rng = $ws.Range("D29:G29")
$rng.NumberFormat = "@"
$ws.Range("G29").Value = '5'
$rng.Style = "Normal"

# Row 30
$rng = $ws.Range("D30:G30")
$rng.NumberFormat = "@"
$ws.Range("G30").Value = '5'
$rng.Style = "Normal"

# Row 31
$rng = $ws.Range("D31:G31")
$rng.NumberFormat = "@"
$ws.Range("G31").Value = '5'
$rng.Style = "Normal"

# Row 32
$rng = $ws.Range("D32:G32")
$rng.NumberFormat = "@"
$ws.Range("G32").Value = '5'
$rng.Style = "Normal"

# Row 33
$rng = $ws.Range("D33:G33")
$rng.NumberFormat = "@"
$ws.Range("G33").Value = '5'
$rng.Style = "Normal"

# Row 34
$rng = $ws.Range("D34:G34")
$rng.NumberFormat = "@"
$ws.Range("G34").Value = '5'
$rng.Style = "Normal"

# Row 35
$rng = $ws.Range("D35:G35")
$rng.NumberFormat = "@"
$ws.Range("G35").Value = '5'
$rng.Style = "Normal"

# Row 36
$rng = $ws.Range("D36:G36")
$rng.NumberFormat = "@"
$ws.Range("G36").Value = '5'
$rng.Style = "Normal"

# Row 37
$rng = $ws.Range("D37:G37")
$rng.NumberFormat = "@"
$ws.Range("G37").Value = '5'
$rng.Style = "Normal"

# Row 38
$rng = $ws.Range("D38:G38")
$rng.NumberFormat = "@"
$ws.Range("D38").Value = '0.02591'
$ws.Range("E38").Value = '-7.39%'
$ws.Range("G38").Value = '5'
$rng.Style = "Normal"

# Row 39
$rng = $ws.Range("D39:G39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = '0.05531'
$ws.Range("E39").Value = '-4.55%'
$ws.Range("G39").Value = '5'
$rng.Style = "Normal"

# Row 40
$rng = $ws.Range("D40:G40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = '0.007818'
$ws.Range("E40").Value = '0.56%'
$ws.Range("G40").Value = '5'
$rng.Style = "Normal"

# Row 41
$rng = $ws.Range("D41:G41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = '0.1393'
$ws.Range("E41").Value = '-2.97%'
$ws.Range("G41").Value = '5'
$rng.Style = "Normal"

# Row 42
$rng = $ws.Range("D42:G42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = '0.006607'
$ws.Range("E42").Value = '-9.27%'
$ws.Range("G42").Value = '5'
$rng.Style = "Normal"

# Row 43
$rng = $ws.Range("D43:G43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = '0.002122'
$ws.Range("E43").Value = '7.38%'
$ws.Range("G43").Value = '5'
$rng.Style = "Normal"

# Row 44
$rng = $ws.Range("D44:G44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = '0.008471'
$ws.Range("E44").Value = '1.36%'
$ws.Range("G44").Value = '5'
$rng.Style = "Normal"

# Row 45
$rng = $ws.Range("D45:G45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = '0.3456'
$ws.Range("E45").Value = '1.79%'
$ws.Range("G45").Value = '5'
$rng.Style = "Normal"

# Row 46
$rng = $ws.Range("D46:G46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = '0.00006974'
$ws.Range("E46").Value = '-5.78%'
$ws.Range("G46").Value = '5'
$rng.Style = "Normal"

# Row 47
$rng = $ws.Range("D47:G47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000753'
$ws.Range("E47").Value = '0.24%'
$ws.Range("G47").Value = '5'
$rng.Style = "Normal"

# Row 48
$rng = $ws.Range("D48:G48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = '0.003477'
$ws.Range("E48").Value = '-0.84%'
$ws.Range("G48").Value = '5'
$rng.Style = "Normal"

# Row 49
$rng = $ws.Range("D49:G49")
$rng.NumberFormat = "@"
$ws.Range("D49").Value = '0.003536'
$ws.Range("E49").Value = '0.93%'
$ws.Range("G49").Value = '5'
$rng.Style = "Normal"

# Row 50
$rng = $ws.Range("D50:G50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = '0.00002108'
$ws.Range("E50").Value = '0.24%'
$ws.Range("G50").Value = '5'
$rng.Style = "Normal"

# Row 51
$rng = $ws.Range("D51:G51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = '0.0002008'
$ws.Range("E51").Value = '0.24%'
$ws.Range("G51").Value = '5'
$rng.Style = "Normal"
